$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1716781.4
$ws.Range("I9").Value = 2199.2222
$ws.Range("J9").Value = 4803029.5
$ws.Range("K9").Value = 2199.2222
$ws.Range("L9").Value = 4803029.5
$ws.Range("M9").Value = -2030.2222
$ws.Range("N9").Value = -4803367.5
$ws.Range("H12").Value = 600
$ws.Range("I12").Value = 600
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -430
$ws.Range("N12").Value = ""
$ws.Range("H17").Value = 2491.6516
$ws.Range("J17").Value = 2491.6516
$ws.Range("L17").Value = 7474.9548
$ws.Range("N17").Value = -7810.9548
$ws.Range("H18").Value = 2788.1667
$ws.Range("J18").Value = 2076.3333
$ws.Range("L18").Value = 2076.3333
$ws.Range("N18").Value = -2644.3333
$ws.Range("H33").Value = 6200.769
$ws.Range("I33").Value = 7255.5454
$ws.Range("K33").Value = 7255.5454
$ws.Range("M33").Value = -7026.5454
$ws.Range("H40").Value = 7719.3335
$ws.Range("I40").Value = 6973.5
$ws.Range("J40").Value = 8316
$ws.Range("K40").Value = 6973.5
$ws.Range("L40").Value = 8316
$ws.Range("M40").Value = -6798.5
$ws.Range("N40").Value = -8666
$ws.Range("H41").Value = 576.9286
$ws.Range("I41").Value = 359.8
$ws.Range("K41").Value = 359.8
$ws.Range("M41").Value = 80.19999999999999
$ws.Range("H43").Value = 10718.235
$ws.Range("I43").Value = 30750
$ws.Range("J43").Value = 8047.3335
$ws.Range("K43").Value = 30750
$ws.Range("L43").Value = 8047.3335
$ws.Range("M43").Value = -30681
$ws.Range("N43").Value = -8185.3335
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H64").Value = 6535.75
$ws.Range("J64").Value = 7809.1665
$ws.Range("L64").Value = 7809.1665
$ws.Range("N64").Value = -8305.166499999999
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H67").Value = 6535.75
$ws.Range("J67").Value = 7809.1665
$ws.Range("L67").Value = 7809.1665
$ws.Range("N67").Value = -9525.166499999999
$ws.Range("H70").Value = 1941145.5
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 2117340.5
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 6352021.5
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -6352561.5
$ws.Range("H73").Value = 1941145.5
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 2117340.5
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 6352021.5
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -6353893.5
$ws.Range("H86").Value = 3815.8333
$ws.Range("J86").Value = 4179
$ws.Range("L86").Value = 4179
$ws.Range("N86").Value = -6425
$ws.Range("H89").Value = 3815.8333
$ws.Range("J89").Value = 4179
$ws.Range("L89").Value = 20895
$ws.Range("N89").Value = -32127
$ws.Range("H99").Value = 298.75
$ws.Range("I99").Value = 298.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 896.25
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 601.75
$ws.Range("N99").Value = ""
$ws.Range("H100").Value = 3816.8462
$ws.Range("J100").Value = 6249.857
$ws.Range("L100").Value = 6249.857
$ws.Range("N100").Value = -7331.857
$ws.Range("H107").Value = 19231862
$ws.Range("I107").Value = 20000656
$ws.Range("K107").Value = 20000656
$ws.Range("M107").Value = -19998736
$ws.Range("H111").Value = 606.0833
$ws.Range("I111").Value = 652.8182
$ws.Range("K111").Value = 1958.4546
$ws.Range("M111").Value = 1108.5454
$ws.Range("H112").Value = 2674.1924
$ws.Range("J112").Value = 2674.1924
$ws.Range("L112").Value = 8022.5772
$ws.Range("N112").Value = -10238.5772
$ws.Range("H125").Value = 4000
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").Value = ""
$ws.Range("H132").Value = 796.6271400000001
$ws.Range("I132").Value = 686.8302
$ws.Range("J132").Value = 1766.5
$ws.Range("K132").Value = 2060.4906
$ws.Range("L132").Value = 5299.5
$ws.Range("M132").Value = 469.5093999999999
$ws.Range("N132").Value = -10359.5
$ws.Range("H135").Value = 958.6111
$ws.Range("I135").Value = 954.05884
$ws.Range("K135").Value = 8586.529560000001
$ws.Range("M135").Value = -6051.529560000001
$ws.Range("H137").Value = 5920.788
$ws.Range("I137").Value = 3611.2917
$ws.Range("K137").Value = 10833.8751
$ws.Range("M137").Value = -8283.875100000001
$ws.Range("H138").Value = 3394.2134
$ws.Range("I138").Value = 3071.0557
$ws.Range("J138").Value = 3692.513
$ws.Range("K138").Value = 9213.167099999999
$ws.Range("L138").Value = 11077.539
$ws.Range("M138").Value = -4073.167099999999
$ws.Range("N138").Value = -21357.539
$ws.Range("H141").Value = 1016.6667
$ws.Range("I141").Value = 1016.6667
$ws.Range("K141").Value = 3050.0001
$ws.Range("M141").Value = 2129.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3540221.2
$ws.Range("I2").Value = 4382269
$ws.Range("K2").Value = 4382269
$ws.Range("M2").Value = -4382156
$ws.Range("H24").Value = 30000
$ws.Range("J24").Value = 30000
$ws.Range("L24").Value = 30000
$ws.Range("N24").Value = -30748
$ws.Range("H32").Value = 3034780.2
$ws.Range("I32").Value = 3129101.5
$ws.Range("K32").Value = 3129101.5
$ws.Range("M32").Value = -3128814.5
$ws.Range("H45").Value = 1665.4
$ws.Range("I45").Value = 1176
$ws.Range("K45").Value = 1176
$ws.Range("M45").Value = -799
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("H61").Value = 3869.7637
$ws.Range("I61").Value = 3192.3137
$ws.Range("K61").Value = 3192.3137
$ws.Range("M61").Value = -2980.3137
$ws.Range("H62").Value = 70249
$ws.Range("J62").Value = 70249
$ws.Range("L62").Value = 70249
$ws.Range("N62").Value = -71497
$ws.Range("H63").Value = 9513.857
$ws.Range("I63").Value = 915.6667
$ws.Range("K63").Value = 915.6667
$ws.Range("M63").Value = -229.6667
$ws.Range("H65").Value = 70249
$ws.Range("J65").Value = 70249
$ws.Range("L65").Value = 210747
$ws.Range("N65").Value = -216987
$ws.Range("H66").Value = 9513.857
$ws.Range("I66").Value = 915.6667
$ws.Range("K66").Value = 4578.3335
$ws.Range("M66").Value = -1146.3335
$ws.Range("H74").Value = 272734.06
$ws.Range("I74").Value = 358557.75
$ws.Range("K74").Value = 358557.75
$ws.Range("M74").Value = -357683.75
$ws.Range("H77").Value = 272734.06
$ws.Range("I77").Value = 358557.75
$ws.Range("K77").Value = 1792788.75
$ws.Range("M77").Value = -1788420.75
$ws.Range("H97").Value = 1267712.6
$ws.Range("I97").Value = 1484892
$ws.Range("J97").Value = 181815.6
$ws.Range("K97").Value = 1484892
$ws.Range("L97").Value = 181815.6
$ws.Range("M97").Value = -1484396
$ws.Range("N97").Value = -182807.6
$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164
$ws.Range("H110").Value = 11365986
$ws.Range("H116").Value = 3540221.2
$ws.Range("I116").Value = 4382269
$ws.Range("K116").Value = 4382269
$ws.Range("M116").Value = -4379975
$ws.Range("H124").Value = 30214.5
$ws.Range("I124").Value = 20000
$ws.Range("J124").Value = 40429
$ws.Range("K124").Value = 20000
$ws.Range("L124").Value = 40429
$ws.Range("M124").Value = -15090
$ws.Range("N124").Value = -50249
$ws.Range("H136").Value = 3869.7637
$ws.Range("I136").Value = 3192.3137
$ws.Range("K136").Value = 9576.9411
$ws.Range("M136").Value = -7026.9411

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3540221.2
$ws.Range("I3").Value = 4382269
$ws.Range("K3").Value = 4382269
$ws.Range("M3").Value = -4382155
$ws.Range("H11").Value = 413.85715
$ws.Range("J11").Value = 400
$ws.Range("L11").Value = 400
$ws.Range("N11").Value = -680
$ws.Range("H57").Value = 103999.664
$ws.Range("J57").Value = 103999.664
$ws.Range("L57").Value = 103999.664
$ws.Range("N57").Value = -105439.664
$ws.Range("H86").Value = 75590.74000000001
$ws.Range("I86").Value = 1528.1875
$ws.Range("J86").Value = 183318.1
$ws.Range("K86").Value = 1528.1875
$ws.Range("L86").Value = 183318.1
$ws.Range("M86").Value = -405.1875
$ws.Range("N86").Value = -185564.1
$ws.Range("H89").Value = 75590.74000000001
$ws.Range("I89").Value = 1528.1875
$ws.Range("J89").Value = 183318.1
$ws.Range("K89").Value = 7640.9375
$ws.Range("L89").Value = 916590.5
$ws.Range("M89").Value = -2024.9375
$ws.Range("N89").Value = -927822.5
$ws.Range("H94").Value = 1139.7142
$ws.Range("I94").Value = 1267.7142
$ws.Range("K94").Value = 1267.7142
$ws.Range("M94").Value = -816.7141999999999
$ws.Range("H99").Value = 8583.333000000001
$ws.Range("I99").Value = 1321.4286
$ws.Range("J99").Value = 34000
$ws.Range("K99").Value = 1321.4286
$ws.Range("L99").Value = 34000
$ws.Range("M99").Value = 176.5714
$ws.Range("N99").Value = -36996
$ws.Range("H132").Value = 107774.5
$ws.Range("J132").Value = 107774.5
$ws.Range("L132").Value = 107774.5
$ws.Range("N132").Value = -117894.5
$ws.Range("H136").Value = 103999.664
$ws.Range("J136").Value = 103999.664
$ws.Range("L136").Value = 103999.664
$ws.Range("N136").Value = -114199.664

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3820.5
$ws.Range("I16").Value = 2825
$ws.Range("K16").Value = 2825
$ws.Range("M16").Value = -2538
$ws.Range("H20").Value = 43333.332
$ws.Range("J20").Value = 43333.332
$ws.Range("L20").Value = 43333.332
$ws.Range("N20").Value = -43805.332
$ws.Range("H22").Value = 3083.7273
$ws.Range("I22").Value = 507
$ws.Range("J22").Value = 4050
$ws.Range("K22").Value = 507
$ws.Range("L22").Value = 4050
$ws.Range("M22").Value = -157
$ws.Range("N22").Value = -4750
$ws.Range("H30").Value = 43333.332
$ws.Range("J30").Value = 43333.332
$ws.Range("L30").Value = 43333.332
$ws.Range("N30").Value = -43515.332
$ws.Range("H31").Value = 23813436
$ws.Range("I31").Value = 55557076
$ws.Range("J31").Value = 5707.4585
$ws.Range("K31").Value = 55557076
$ws.Range("L31").Value = 5707.4585
$ws.Range("M31").Value = -55556781
$ws.Range("N31").Value = -6297.4585
$ws.Range("H34").Value = 23813436
$ws.Range("I34").Value = 55557076
$ws.Range("J34").Value = 5707.4585
$ws.Range("K34").Value = 55557076
$ws.Range("L34").Value = 5707.4585
$ws.Range("M34").Value = -55556874
$ws.Range("N34").Value = -6111.4585
$ws.Range("H43").Value = 98571.28999999999
$ws.Range("J43").Value = 98571.28999999999
$ws.Range("L43").Value = 98571.28999999999
$ws.Range("N43").Value = -98939.28999999999
$ws.Range("H58").Value = 4550.7144
$ws.Range("I58").Value = 4573.724
$ws.Range("K58").Value = 4573.724
$ws.Range("M58").Value = -4370.724
$ws.Range("H62").Value = 12253.772
$ws.Range("J62").Value = 15506.1875
$ws.Range("L62").Value = 15506.1875
$ws.Range("N62").Value = -16754.1875
$ws.Range("H65").Value = 12253.772
$ws.Range("J65").Value = 15506.1875
$ws.Range("L65").Value = 77530.9375
$ws.Range("N65").Value = -83770.9375
$ws.Range("H101").Value = 98571.28999999999
$ws.Range("J101").Value = 98571.28999999999
$ws.Range("L101").Value = 98571.28999999999
$ws.Range("N101").Value = -105061.29
$ws.Range("H113").Value = 3820.5
$ws.Range("I113").Value = 2825
$ws.Range("K113").Value = 2825
$ws.Range("M113").Value = -655
$ws.Range("H123").Value = 54922.332
$ws.Range("J123").Value = 54922.332
$ws.Range("L123").Value = 54922.332
$ws.Range("N123").Value = -64722.332
$ws.Range("H125").Value = 101000
$ws.Range("J125").Value = 101000
$ws.Range("L125").Value = 101000
$ws.Range("N125").Value = -105920
$ws.Range("H128").Value = 43333.332
$ws.Range("J128").Value = 43333.332
$ws.Range("L128").Value = 43333.332
$ws.Range("N128").Value = -53293.332
$ws.Range("H132").Value = 5315.6
$ws.Range("I132").Value = 4538.294
$ws.Range("J132").Value = 6967.375
$ws.Range("K132").Value = 13614.882
$ws.Range("L132").Value = 20902.125
$ws.Range("M132").Value = -11084.882
$ws.Range("N132").Value = -25962.125
$ws.Range("H134").Value = 5719.657
$ws.Range("I134").Value = 5764.0347
$ws.Range("J134").Value = 5505.1665
$ws.Range("K134").Value = 17292.1041
$ws.Range("L134").Value = 16515.4995
$ws.Range("M134").Value = -14757.1041
$ws.Range("N134").Value = -21585.4995
$ws.Range("H136").Value = 4550.7144
$ws.Range("I136").Value = 4573.724
$ws.Range("K136").Value = 13721.172
$ws.Range("M136").Value = -11171.172

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 787.125
$ws.Range("I5").Value = 685.4286
$ws.Range("K5").Value = 2056.2858
$ws.Range("M5").Value = -1944.2858
$ws.Range("H75").Value = 144.5
$ws.Range("J75").Value = 144.5
$ws.Range("L75").Value = 433.5
$ws.Range("N75").Value = -2429.5
$ws.Range("H78").Value = 144.5
$ws.Range("J78").Value = 144.5
$ws.Range("L78").Value = 1300.5
$ws.Range("N78").Value = -11284.5
$ws.Range("H88").Value = 3523
$ws.Range("J88").Value = 9999
$ws.Range("L88").Value = 29997
$ws.Range("N88").Value = -30853
$ws.Range("H91").Value = 3523
$ws.Range("J91").Value = 9999
$ws.Range("L91").Value = 29997
$ws.Range("N91").Value = -32961
$ws.Range("H114").Value = 599.6875
$ws.Range("I114").Value = 395.5
$ws.Range("K114").Value = 1186.5
$ws.Range("M114").Value = 2067.5
$ws.Range("H122").Value = 142867500
$ws.Range("I122").Value = 250017580
$ws.Range("K122").Value = 2250158220
$ws.Range("M122").Value = -2250155770
$ws.Range("H135").Value = 787.125
$ws.Range("I135").Value = 685.4286
$ws.Range("K135").Value = 6168.8574
$ws.Range("M135").Value = -3633.8574

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""
$ws.Range("H11").Value = 4099875.2
$ws.Range("I11").Value = 555566.9
$ws.Range("J11").Value = 10479630
$ws.Range("K11").Value = 555566.9
$ws.Range("L11").Value = 10479630
$ws.Range("M11").Value = -555427.9
$ws.Range("N11").Value = -10479908
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10346
$ws.Range("H30").Value = 10000
$ws.Range("J30").Value = 10000
$ws.Range("L30").Value = 10000
$ws.Range("N30").Value = -10210
$ws.Range("H70").Value = 6537.846
$ws.Range("I70").Value = 5566.222
$ws.Range("K70").Value = 5566.222
$ws.Range("M70").Value = -5296.222
$ws.Range("H73").Value = 6537.846
$ws.Range("I73").Value = 5566.222
$ws.Range("K73").Value = 5566.222
$ws.Range("M73").Value = -4630.222
$ws.Range("H97").Value = 1863.1395
$ws.Range("I97").Value = 1671.6285
$ws.Range("J97").Value = 2701
$ws.Range("K97").Value = 1671.6285
$ws.Range("L97").Value = 2701
$ws.Range("M97").Value = -1175.6285
$ws.Range("N97").Value = -3693
$ws.Range("H101").Value = 21198.2
$ws.Range("J101").Value = 21198.2
$ws.Range("L101").Value = 21198.2
$ws.Range("N101").Value = -27688.2
$ws.Range("H102").Value = 1342.4445
$ws.Range("I102").Value = 1239.7858
$ws.Range("K102").Value = 1239.7858
$ws.Range("M102").Value = 382.2141999999999
$ws.Range("H122").Value = 2716.9092
$ws.Range("I122").Value = 2814.7097
$ws.Range("J122").Value = 2483.6924
$ws.Range("K122").Value = 8444.1291
$ws.Range("L122").Value = 7451.0772
$ws.Range("M122").Value = -5994.1291
$ws.Range("N122").Value = -12351.0772
$ws.Range("H123").Value = 42500
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -84900
$ws.Range("H132").Value = 4979.727
$ws.Range("I132").Value = 2918
$ws.Range("K132").Value = 8754
$ws.Range("M132").Value = -6224

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 18000
$ws.Range("I23").Value = 18000
$ws.Range("K23").Value = 18000
$ws.Range("M23").Value = -17770
$ws.Range("H46").Value = 4626.1836
$ws.Range("I46").Value = 1391
$ws.Range("J46").Value = 5562.684
$ws.Range("K46").Value = 1391
$ws.Range("L46").Value = 5562.684
$ws.Range("M46").Value = -1203
$ws.Range("N46").Value = -5938.684
$ws.Range("H61").Value = 1794
$ws.Range("I61").Value = 2013.7778
$ws.Range("K61").Value = 2013.7778
$ws.Range("M61").Value = -1811.7778
$ws.Range("H68").Value = 5824.8
$ws.Range("I68").Value = 3062.5
$ws.Range("K68").Value = 3062.5
$ws.Range("M68").Value = -2313.5
$ws.Range("H71").Value = 5824.8
$ws.Range("I71").Value = 3062.5
$ws.Range("K71").Value = 15312.5
$ws.Range("M71").Value = -11568.5
$ws.Range("H93").Value = 3150
$ws.Range("J93").Value = 4000
$ws.Range("L93").Value = 4000
$ws.Range("N93").Value = -6496
$ws.Range("H113").Value = 1794
$ws.Range("I113").Value = 2013.7778
$ws.Range("K113").Value = 2013.7778
$ws.Range("M113").Value = 156.2221999999999
$ws.Range("H122").Value = 38464664
$ws.Range("I122").Value = 50002700
$ws.Range("K122").Value = 150008100
$ws.Range("M122").Value = -150005650
$ws.Range("H132").Value = 6738.2173
$ws.Range("I132").Value = 1865.5
$ws.Range("K132").Value = 5596.5
$ws.Range("M132").Value = -3066.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""
$ws.Range("H19").Value = 14000
$ws.Range("J19").Value = 14000
$ws.Range("L19").Value = 14000
$ws.Range("N19").Value = -14348
$ws.Range("H46").Value = 86475
$ws.Range("J46").Value = 86475
$ws.Range("L46").Value = 86475
$ws.Range("N46").Value = -86937
$ws.Range("H62").Value = 12145.733
$ws.Range("I62").Value = 13881.833
$ws.Range("J62").Value = 10988.333
$ws.Range("K62").Value = 13881.833
$ws.Range("L62").Value = 10988.333
$ws.Range("M62").Value = -13257.833
$ws.Range("N62").Value = -12236.333
$ws.Range("H65").Value = 12145.733
$ws.Range("I65").Value = 13881.833
$ws.Range("J65").Value = 10988.333
$ws.Range("K65").Value = 69409.16500000001
$ws.Range("L65").Value = 54941.665
$ws.Range("M65").Value = -66289.16500000001
$ws.Range("N65").Value = -61181.665
$ws.Range("H81").Value = 11331.765
$ws.Range("I81").Value = 3336.7334
$ws.Range("J81").Value = 14663.027
$ws.Range("K81").Value = 6673.4668
$ws.Range("L81").Value = 29326.054
$ws.Range("M81").Value = -5612.4668
$ws.Range("N81").Value = -31448.054
$ws.Range("H84").Value = 11331.765
$ws.Range("I84").Value = 3336.7334
$ws.Range("J84").Value = 14663.027
$ws.Range("K84").Value = 33367.334
$ws.Range("L84").Value = 146630.27
$ws.Range("M84").Value = -28063.334
$ws.Range("N84").Value = -157238.27
$ws.Range("H96").Value = 6231.077
$ws.Range("J96").Value = 7143.143
$ws.Range("L96").Value = 7143.143
$ws.Range("N96").Value = -9889.143
$ws.Range("H107").Value = 1077.125
$ws.Range("J107").Value = 2975.4
$ws.Range("L107").Value = 8926.200000000001
$ws.Range("N107").Value = -12766.2
$ws.Range("H122").Value = 3971.2856
$ws.Range("I122").Value = 3971.2856
$ws.Range("K122").Value = 11913.8568
$ws.Range("M122").Value = -9463.856800000001
$ws.Range("H132").Value = 5158.0967
$ws.Range("I132").Value = 4583.4443
$ws.Range("K132").Value = 13750.3329
$ws.Range("M132").Value = -11220.3329
$ws.Range("H134").Value = 86475
$ws.Range("J134").Value = 86475
$ws.Range("L134").Value = 259425
$ws.Range("N134").Value = -264495
$ws.Range("H136").Value = 2751.4
$ws.Range("I136").Value = 997.5833
$ws.Range("J136").Value = 9766.666999999999
$ws.Range("K136").Value = 2992.7499
$ws.Range("L136").Value = 29300.001
$ws.Range("M136").Value = -442.7498999999998
$ws.Range("N136").Value = -34400.001
